$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values (column D),
# matching the original workbook's text representation.
$textForceCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D16", "D17", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D30", "D31", "D35", "D36", "D37", "D38", "D41", "D42", "D45", "D46", "D49", "D51")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated price (D) and volume-change (E) values
$ws.Range("D2").Value = "41.208.68"
$ws.Range("E2").Value = "  -5.80%  "
$ws.Range("D3").Value = "2.220.23"
$ws.Range("E3").Value = "  -5.59%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "245.63"
$ws.Range("E5").Value = "  +2.65%  "
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  -5.98%  "
$ws.Range("D7").Value = "70.27"
$ws.Range("E7").Value = "  -5.04%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  -6.82%  "
$ws.Range("D10").Value = "38.94"
$ws.Range("E10").Value = "  +4.77%  "
$ws.Range("D11").Value = "0.0951"
$ws.Range("E11").Value = "  -6.82%  "
$ws.Range("D12").Value = "58.35"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("D14").Value = "6.76"
$ws.Range("E14").Value = "  -7.49%  "
$ws.Range("D15").Value = "2.547.24"
$ws.Range("E15").Value = "  -5.81%  "
$ws.Range("D16").Value = "14.84"
$ws.Range("E16").Value = "  -9.22%  "
$ws.Range("D17").Value = "0.843"
$ws.Range("E17").Value = "  -8.56%  "
$ws.Range("D18").Value = "2.221.89"
$ws.Range("E18").Value = "  -5.96%  "
$ws.Range("D19").Value = "41.269.48"
$ws.Range("E19").Value = "  -5.48%  "
$ws.Range("E20").Value = "  -7.93%  "
$ws.Range("D21").Value = "72.51"
$ws.Range("E21").Value = "  -5.85%  "
$ws.Range("D22").Value = "6.08"
$ws.Range("E22").Value = "  -7.69%  "
$ws.Range("D23").Value = "231.95"
$ws.Range("E23").Value = "  -8.31%  "
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +11.89%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "3.68"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "2.43"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("D28").Value = "9.87"
$ws.Range("E28").Value = "  -7.05%  "
$ws.Range("E29").Value = "  -4.99%  "
$ws.Range("D30").Value = "172.25"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").Value = "20.53"
$ws.Range("E31").Value = "  -7.76%  "
$ws.Range("E32").Value = "  -7.48%  "
$ws.Range("E33").Value = "  -6.91%  "
$ws.Range("E34").Value = "  -5.38%  "
$ws.Range("D35").Value = "5.25"
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("D36").Value = "4.61"
$ws.Range("E36").Value = "  -10.05%  "
$ws.Range("D37").Value = "3.90"
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("D38").Value = "24.30"
$ws.Range("E38").Value = "  +17.13%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  -5.18%  "
$ws.Range("D41").Value = "5.86"
$ws.Range("E41").Value = "  -10.92%  "
$ws.Range("D42").Value = "65.69"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("E43").Value = "  -9.29%  "
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("D45").Value = "8.83"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").Value = "10.91"
$ws.Range("E46").Value = "  +11.62%  "
$ws.Range("E47").Value = "  -6.41%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Value = "4.55"
$ws.Range("E49").Value = "  +5.22%  "
$ws.Range("E50").Value = "  -5.47%  "
$ws.Range("D51").Value = "1.10"
$ws.Range("E51").Value = "  -4.94%  "
